$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the three runs of the "Transferir todo o cadastro de
# clientes..." bullet into a single run (same resulting text/formatting).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "do sistema da Zé Pequeno eletro para o CRM principal da VendasGeral, e depreciar o atual sistema de clientes da Zé Pequeno;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "do sistema da Zé Pequeno eletro para o CRM principal da VendasGeral, e depreciar o atual sistema de clientes da Zé Pequeno;",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Análise de Gaps" section rewrite.
# ---------------------------------------------------------------------------

# Donor range whose run carries the rPr (rFonts/color/kern/sz/szCs/lang) we
# need to stamp onto the newly-created runs/paragraphs below.
$donor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("- Cancelar o")) {
        $donor = $p.Range
        break
    }
}

# Step 2a: remove the stray "As princi" paragraph entirely (paragraph mark
# included), leaving the following paragraph's pPr untouched.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("As princi")) {
        $p.Range.Delete()
        break
    }
}

# Step 2b: split the "Demonstrar as principais..." paragraph into four
# paragraphs with new content, reusing the original paragraph's pPr
# (LOnormal / firstLine 720 / jc both) for each of them.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Demonstrar as principais")) {
        $p.Range.Text = "As principais diferen`r" + `
            "PARA2`r" + `
            "PARA3`r" + `
            "PARA4"

        # --- paragraph 1: second run "ças entre a arquitetura atual e futura são:" ---
        $tailStart = $p.Range.End - 1
        $r = $d.Range($tailStart, $tailStart)
        $r.InsertAfter("ças entre a arquitetura atual e futura são:")
        $tailEnd = $p.Range.End - 1
        $r = $d.Range($tailStart, $tailEnd)
        $r.FormattedText = $donor.FormattedText
        $tailEnd2 = $p.Range.End - 1
        $r = $d.Range($tailStart, $tailEnd2)
        $r.Text = "ças entre a arquitetura atual e futura são:"

        $p2 = $p.Next()
        $p2.Range.Text = $p2.Range.Text  # no-op, keep reference fresh
        $r = $p2.Range
        $r.FormattedText = $donor.FormattedText
        $r = $p2.Range
        $r.Text = "- Diminuição no número de aplicações, que resulta em uma redução do custo de gerencimanto das mesmas e num maior reuso e escalabidade;"

        $p3 = $p2.Next()
        $r = $p3.Range
        $r.FormattedText = $donor.FormattedText
        $r = $p3.Range
        $r.Text = "- Foco na lucratividade, já que produtos e cupons agora serão até certo ponto tratados como a mesma coisa pelos processos do neǵocio;"

        $p4 = $p3.Next()
        $r = $p4.Range
        $r.FormattedText = $donor.FormattedText
        $r = $p4.Range
        $r.Text = "- Expansão do mercado consumidor, utilizando-se de artificios de marketing já existente para capturar clientes de produtos para cupons e vice versa;"

        break
    }
}
